$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to store the value as text (matches the source
    # data which always uses inline/shared text strings, even for
    # numeric-looking prices), then restore the default "Normal"
    # style so no stray number-format styling is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '63.879.48'
Set-TextValue $ws.Range('E2') '  -1.55%  '
Set-TextValue $ws.Range('D3') '3.061.81'
Set-TextValue $ws.Range('E3') '  -1.21%  '
Set-TextValue $ws.Range('E4') '  +0.00%  '
Set-TextValue $ws.Range('D5') '558.63'
Set-TextValue $ws.Range('E5') '  -0.38%  '
Set-TextValue $ws.Range('D6') '142.46'
Set-TextValue $ws.Range('E6') '  -1.22%  '
Set-TextValue $ws.Range('E7') '  +0.01%  '
Set-TextValue $ws.Range('D8') '3.059.16'
Set-TextValue $ws.Range('E8') '  -1.18%  '
Set-TextValue $ws.Range('E9') '  +3.33%  '
Set-TextValue $ws.Range('E10') '  +0.48%  '
Set-TextValue $ws.Range('D11') '6.19'
Set-TextValue $ws.Range('E11') '  -3.57%  '
Set-TextValue $ws.Range('E12') '  +1.72%  '
Set-TextValue $ws.Range('E13') '  +1.12%  '
Set-TextValue $ws.Range('D14') '35.26'
Set-TextValue $ws.Range('E14') '  -0.24%  '
Set-TextValue $ws.Range('D15') '3.563.19'
Set-TextValue $ws.Range('E15') '  -1.06%  '
Set-TextValue $ws.Range('D16') '63.928.29'
Set-TextValue $ws.Range('E16') '  -1.49%  '
Set-TextValue $ws.Range('D17') '3.061.25'
Set-TextValue $ws.Range('E17') '  -1.26%  '
Set-TextValue $ws.Range('E18') '  +0.09%  '
Set-TextValue $ws.Range('D19') '6.78'
Set-TextValue $ws.Range('E19') '  +0.07%  '
Set-TextValue $ws.Range('D20') '486.77'
Set-TextValue $ws.Range('E20') '  +1.36%  '
Set-TextValue $ws.Range('D21') '14.36'
Set-TextValue $ws.Range('E21') '  +3.73%  '
Set-TextValue $ws.Range('E22') '  +0.06%  '
Set-TextValue $ws.Range('D23') '14.60'
Set-TextValue $ws.Range('E23') '  +7.95%  '
Set-TextValue $ws.Range('D24') '7.54'
Set-TextValue $ws.Range('E24') '  -0.44%  '
Set-TextValue $ws.Range('D25') '82.70'
Set-TextValue $ws.Range('E25') '  +1.86%  '
Set-TextValue $ws.Range('E26') '  +0.09%  '
Set-TextValue $ws.Range('D27') '2.80'
Set-TextValue $ws.Range('E27') '  +0.10%  '
Set-TextValue $ws.Range('D28') '8.14'
Set-TextValue $ws.Range('E28') '  -0.57%  '
Set-TextValue $ws.Range('E29') '  -0.80%  '
Set-TextValue $ws.Range('E30') '  +0.01%  '
Set-TextValue $ws.Range('D31') '26.51'
Set-TextValue $ws.Range('E31') '  +1.44%  '
Set-TextValue $ws.Range('D32') '1.16'
Set-TextValue $ws.Range('E32') '  +0.57%  '
Set-TextValue $ws.Range('D33') '2.52'
Set-TextValue $ws.Range('E33') '  +0.49%  '
Set-TextValue $ws.Range('E34') '  +1.17%  '
Set-TextValue $ws.Range('E35') '  +1.11%  '
Set-TextValue $ws.Range('D36') '54.93'
Set-TextValue $ws.Range('E36') '  -0.19%  '
Set-TextValue $ws.Range('E37') '  +0.86%  '
Set-TextValue $ws.Range('D38') '444.59'
Set-TextValue $ws.Range('E38') '  -5.54%  '
Set-TextValue $ws.Range('D39') '0.0814'
Set-TextValue $ws.Range('E39') '  -2.66%  '
Set-TextValue $ws.Range('B40') 'Maker'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D40') '3.025.34'
Set-TextValue $ws.Range('E40') '  +1.44%  '
Set-TextValue $ws.Range('B41') 'dogwifhat'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D41') '2.78'
Set-TextValue $ws.Range('E41') '  -5.70%  '
Set-TextValue $ws.Range('E42') '  +0.95%  '
Set-TextValue $ws.Range('E43') '  +1.37%  '
Set-TextValue $ws.Range('D44') '0.275'
Set-TextValue $ws.Range('E44') '  +5.89%  '
Set-TextValue $ws.Range('D45') '2.26'
Set-TextValue $ws.Range('D46') '27.66'
Set-TextValue $ws.Range('E46') '  -2.26%  '
Set-TextValue $ws.Range('D48') '0.113'
Set-TextValue $ws.Range('E48') '  +0.97%  '
Set-TextValue $ws.Range('E49') '  -2.53%  '
Set-TextValue $ws.Range('D50') '117.99'
Set-TextValue $ws.Range('E50') '  +0.21%  '
Set-TextValue $ws.Range('E51') '  +2.35%  '
